$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("https://gateway.research.unimelb.edu.au/funding-contracts-and-ethics/funding-and-grants/research-funding-schemes/pivot-rp-funding-database", "https://sites.research.unimelb.edu.au/research-funding/domestic", "Pivot-RP"),
    @("https://sites.research.unimelb.edu.au/research-funding/assets/interim-private-files/2024-CCV-Grants-in-Aid-EOI-online-form-template-for-internal-use.docx", "https://sites.research.unimelb.edu.au/research-funding/domestic/cancer-council-victoria-grants-in-aid", "internal UoM template"),
    @("https://unimelbcloud.sharepoint.com/:v:/t/APPIIARCteam/ESu9uHKuH7JCvYET0LcGETEB_O2RT0ZKbe1JDptiq99KDg?e=UMmzmk", "https://sites.research.unimelb.edu.au/research-funding/domestic/westpac", "here "),
    @("https://gateway.research.unimelb.edu.au/resources/contracts/preparing/costing-pricing-support?SQ_VARIATION_3200695=0", "https://sites.research.unimelb.edu.au/research-funding/domestic/westpac", "Research Costing and Pricing Tool"),
    @("https://gateway.research.unimelb.edu.au/resources/funding/grant-resources/fellowship-salary-gap-management-process", "https://sites.research.unimelb.edu.au/research-funding/domestic/westpac", "Fellowship funding policy"),
    @("https://sites.research.unimelb.edu.au/research-funding/assets/interim-private-files/2023-ECRF-MCRF-HeadofResearchTeam-Template.docx", "https://sites.research.unimelb.edu.au/research-funding/domestic/vca-grants", "ECRF"),
    @("https://gateway.research.unimelb.edu.au/resources/systems-reporting/themis-guides/grants,-contracts-and-submissions/create-a-grant-submission-record", "https://sites.research.unimelb.edu.au/research-funding/domestic/westpac", "Themis submission"),
    @("https://sites.research.unimelb.edu.au/research-funding/assets/interim-private-files/2025_Westpac_Eligibility-Exemption-Request-Calculator_V2.xlsx", "https://sites.research.unimelb.edu.au/research-funding/domestic/westpac", "Eligibility Exemption Calculator"),
    @("https://sites.research.unimelb.edu.au/research-funding/assets/interim-private-files/UoM-Co-Investment-Request-Form_2024-VDOH-VMRAF_updated-21-Nov-23.docx", "https://sites.research.unimelb.edu.au/research-funding/domestic/djpr-victorian-medical-research-acceleration-fund-vmraf-and-mrna-vmraf", "Cash Co-Investment Request Form"),
    @("https://sites.research.unimelb.edu.au/research-funding/assets/interim-private-files/2024-VMRAF-R7-DVCR-Lead-Letter-of-Support-template.docx", "https://sites.research.unimelb.edu.au/research-funding/domestic/djpr-victorian-medical-research-acceleration-fund-vmraf-and-mrna-vmraf", "DVCR Letter of Support template"),
    @("https://gateway.research.unimelb.edu.au/funding-contracts-and-ethics/contracts-and-agreements/prepare-a-contract/research-costing-and-pricing-support", "https://sites.research.unimelb.edu.au/research-funding/domestic/djpr-victorian-medical-research-acceleration-fund-vmraf-and-mrna-vmraf", "Research Costing and Pricing Tool"),
    @("https://sites.research.unimelb.edu.au/research-funding/assets/interim-private-files/UoM-Co-Investment-Request-Form_2025-DJSIR-mRNA-VRAF_updated-1-Oct-24.docx", "https://sites.research.unimelb.edu.au/research-funding/domestic/mrna-vraf", "Cash Co-Investment Request Form"),
    @("https://sites.research.unimelb.edu.au/research-funding/assets/interim-private-files/2025-mRNA-VRAF-R4-DVCR-Lead-Letter-of-Support-template.docx", "https://sites.research.unimelb.edu.au/research-funding/domestic/mrna-vraf", "DVCR Letter of Support template")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
